$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.031031
$ws.Range("H2").Value = 12.093093
$ws.Range("I2").Value = 0.380357182622003
$ws.Range("J2").Value = 0.380357182622003
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.445638539112
$ws.Range("R2").Value = 4.010746852008
$ws.Range("S2").Value = 0.004283771102269286
$ws.Range("T2").Value = 0.004283771102269286
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.031031
$ws.Range("H3").Value = 12.093093
$ws.Range("I3").Value = 0.380357182622003
$ws.Range("J3").Value = 0.380357182622003
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 36.32793623152399
$ws.Range("R3").Value = 326.951426083716
$ws.Range("S3").Value = 0.3492080459283899
$ws.Range("T3").Value = 0.3492080459283899
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.031031
$ws.Range("H4").Value = 12.093093
$ws.Range("I4").Value = 0.380357182622003
$ws.Range("J4").Value = 0.380357182622003
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 2.794790381888999
$ws.Range("R4").Value = 25.153113437001
$ws.Range("S4").Value = 0.02686536559134379
$ws.Range("T4").Value = 0.0268653655913438
$ws.Range("I5").Value = 0.4810839099297969
$ws.Range("J5").Value = 0.4810839099297969
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 0.5636531676186667
$ws.Range("R5").Value = 5.072878508568
$ws.Range("S5").Value = 0.005418205427113096
$ws.Range("T5").Value = 0.005418205427113097
$ws.Range("I6").Value = 0.4810839099297969
$ws.Range("J6").Value = 0.4810839099297969
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.4416858147809183
$ws.Range("T6").Value = 0.4416858147809183
$ws.Range("I7").Value = 0.4810839099297969
$ws.Range("J7").Value = 0.4810839099297969
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 3.534910725452333
$ws.Range("R7").Value = 31.814196529071
$ws.Range("S7").Value = 0.03397988972176555
$ws.Range("T7").Value = 0.03397988972176556
$ws.Range("G8").Value = 1.468449333333333
$ws.Range("H8").Value = 4.405348
$ws.Range("I8").Value = 0.1385589074482
$ws.Range("J8").Value = 0.1385589074482
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 0.1623400106986667
$ws.Range("R8").Value = 1.461060096288
$ws.Range("S8").Value = 0.00156051908786609
$ws.Range("T8").Value = 0.00156051908786609
$ws.Range("G9").Value = 1.468449333333333
$ws.Range("H9").Value = 4.405348
$ws.Range("I9").Value = 0.1385589074482
$ws.Range("J9").Value = 0.1385589074482
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 13.23376916241955
$ws.Range("R9").Value = 119.103922461776
$ws.Range("S9").Value = 0.1272117039631251
$ws.Range("T9").Value = 0.1272117039631251
$ws.Range("G10").Value = 1.468449333333333
$ws.Range("H10").Value = 4.405348
$ws.Range("I10").Value = 0.1385589074482
$ws.Range("J10").Value = 0.1385589074482
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 1.018103823337333
$ws.Range("R10").Value = 9.162934410036
$ws.Range("S10").Value = 0.009786684397208818
$ws.Range("T10").Value = 0.009786684397208818
